$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "pt_max" column (F) is being removed entirely; Excel shifts every
# later column one slot to the left (G->F, H->G, ... M->L) and updates
# the formulas / shared-string references along the way.
[void]$ws.Columns("F:F").Delete()

# The new header row (row 1, now A1:L1) is bold + center-aligned.
$headerRange = $ws.Range("A1:L1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108

# Leave the selection where the author's last save left it.
[void]$ws.Range("H15").Select()
